# Update the two-digit-division answer table: each cell's text is replaced
# in place (same run/paragraph/cell formatting, only the w:t content changes).
$d = $word.ActiveDocument

$d.Content.Find.Execute("52÷7=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=18, 1", 2) | Out-Null
$d.Content.Find.Execute("88÷8=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=21, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷5=8, 1", 2) | Out-Null
$d.Content.Find.Execute("56÷4=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "51÷7=7, 2", 2) | Out-Null
$d.Content.Find.Execute("23÷3=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "39÷6=6, 3", 2) | Out-Null
$d.Content.Find.Execute("62÷2=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷2=44, 1", 2) | Out-Null
$d.Content.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=47, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷2=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=6, 0", 2) | Out-Null
$d.Content.Find.Execute("49÷2=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷6=5, 4", 2) | Out-Null
$d.Content.Find.Execute("25÷2=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "32÷4=8, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷6=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=3, 3", 2) | Out-Null
$d.Content.Find.Execute("37÷4=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷6=6, 1", 2) | Out-Null
$d.Content.Find.Execute("22÷8=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "13÷7=1, 6", 2) | Out-Null
$d.Content.Find.Execute("97÷2=48, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷6=13, 3", 2) | Out-Null
$d.Content.Find.Execute("57÷8=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷4=20, 1", 2) | Out-Null
$d.Content.Find.Execute("72÷7=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2) | Out-Null
$d.Content.Find.Execute("36÷8=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "74÷3=24, 2", 2) | Out-Null
$d.Content.Find.Execute("35÷7=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=6, 6", 2) | Out-Null
$d.Content.Find.Execute("14÷8=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "78÷7=11, 1", 2) | Out-Null
$d.Content.Find.Execute("73÷9=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("36÷3=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷7=11, 1", 2) | Out-Null
$d.Content.Find.Execute("92÷9=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "60÷7=8, 4", 2) | Out-Null
$d.Content.Find.Execute("41÷7=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=42, 0", 2) | Out-Null
$d.Content.Find.Execute("92÷3=30, 2", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=3, 3", 2) | Out-Null
